# Weekly update: insert 3 new daily-price rows for Pimiento (Vega Modelo de
# Temuco) at the top of the data block, pushing the existing rows down.
#
# Before: data rows occupy 1084-1153 (dimension A1:R1153)
# After:  3 new rows are inserted at 1084:1086 (existing rows shift to
#         1087:1156), dimension becomes A1:R1156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current first data row of the block; Excel
# copies formatting (incl. the date style on column D) from the row above.
$ws.Rows("1084:1086").Insert()

# Columns that stay constant across this entire market/product block.
$A = 10
$B = "Vega Modelo de Temuco"
$C = "La Araucanía"
$E = 9
$F = 100112002
$G = "Pimiento"
$R = "Hortaliza"

function Set-DataRow {
    param($row, $D, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q)

    $ws.Range("A$row").Value = $A
    $ws.Range("B$row").Value = $B
    $ws.Range("C$row").Value = $C
    $ws.Range("D$row").Value = $D
    $ws.Range("E$row").Value = $E
    $ws.Range("F$row").Value = $F
    $ws.Range("G$row").Value = $G
    $ws.Range("H$row").Value = $H
    $ws.Range("I$row").Value = $I
    $ws.Range("J$row").Value = $J
    $ws.Range("K$row").Value = $K
    $ws.Range("L$row").Value = $L
    $ws.Range("M$row").Value = $M
    $ws.Range("N$row").Value = $N
    $ws.Range("O$row").Value = $O
    $ws.Range("P$row").Value = $P
    $ws.Range("Q$row").Value = $Q
    $ws.Range("R$row").Value = $R
}

Set-DataRow 1084 44585 "Cuatro cascos verde" "Primera" 235 10000 10000 10000 "`$/caja 15 kilos" "Región del Maule" 667 15
Set-DataRow 1085 44585 "Zafiro rojo" "Primera" 295 23000 25000 23847 "`$/caja 15 kilos" "Región de Arica y Parinacota" 1590 15
Set-DataRow 1086 44585 "Zafiro verde" "Primera" 210 15000 15000 15000 "`$/caja 15 kilos" "Región de Arica y Parinacota" 1000 15
